$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("04/27/2021 18:41:16", "16.344"),
    @("04/27/2021 18:43:06", "16.841"),
    @("04/27/2021 18:51:30", "13.558"),
    @("04/27/2021 18:54:44", "12.788"),
    @("04/27/2021 18:57:56", "12.66")
)

$startRow = 265
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $cellA = $ws.Cells.Item($row, 1)
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellA.Value = $data[$i][0]
    $cellB.Value = $data[$i][1]
}
